$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F (old F "Volume" shifts to G) and a new row
# before old row 6 (old row 6 "Volume" shifts to row 7). This makes room
# for the new "Adj Close" row/column.
$ws.Columns.Item(6).Insert()
$ws.Rows.Item(6).Insert()

# --- Header row (row 1) ---
$ws.Range("F1").Value = "Adj Close"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# --- Row label column (column A) ---
$ws.Range("A6").Value = "Adj Close"
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

# --- Correlation values ---
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.9994160075792081
$ws.Range("D2").Value = 0.9988591541715838
$ws.Range("E2").Value = 0.9985406712857549
$ws.Range("F2").Value = 0.9985406712857549
$ws.Range("G2").Value = 0.6642954180740478

$ws.Range("B3").Value = 0.9994160075792081
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.9988008175769721
$ws.Range("E3").Value = 0.99936633736531
$ws.Range("F3").Value = 0.99936633736531
$ws.Range("G3").Value = 0.6685897361004693

$ws.Range("B4").Value = 0.9988591541715838
$ws.Range("C4").Value = 0.9988008175769721
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.999244360987668
$ws.Range("F4").Value = 0.999244360987668
$ws.Range("G4").Value = 0.6547527316865779

$ws.Range("B5").Value = 0.9985406712857549
$ws.Range("C5").Value = 0.99936633736531
$ws.Range("D5").Value = 0.999244360987668
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.662866677809267

$ws.Range("B6").Value = 0.9985406712857549
$ws.Range("C6").Value = 0.99936633736531
$ws.Range("D6").Value = 0.999244360987668
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.662866677809267

$ws.Range("B7").Value = 0.6642954180740478
$ws.Range("C7").Value = 0.6685897361004693
$ws.Range("D7").Value = 0.6547527316865779
$ws.Range("E7").Value = 0.662866677809267
$ws.Range("F7").Value = 0.662866677809267
$ws.Range("G7").Value = 1
